$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($i = 1; $i -le 50; $i++) {
    $row = $i + 1
    $ws.Cells.Item($row, 5).Value = "Some text for bio $i"
}
